# update template laba rugi
#
# Applies the profit-loss template edit: updated "Report Formula" /
# "Report Operator" strings on the LABA BRUTO (row 11), TOTAL BEBAN
# (row 36), LABA USAHA (row 37) and LABA USAHA SETELAH PAJAK (row 41)
# lines, widens the Report Operator column, bumps a couple of
# "Report Type" values, resizes most data rows, and moves the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (LABA BRUTO): formula/operator now folds in row 1 (PENJUALAN) ---
$ws.Range("I11").Value = "4#5#6#7#8#1"
$ws.Range("J11").Value = '"+#+#+#+#-#-'

# --- Row 36 (TOTAL BEBAN): drop row 19/39, renumber the chain, reflow operators ---
$ws.Range("I36").Value = "14#15#16#17#18#21#22#23#24#25#26#27#28#29#30#31#32#33#34"
$ws.Range("J36").Value = '"+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+'

# --- Row 37 (LABA USAHA): now the full chain through TOTAL BEBAN, Report Type bumped ---
$ws.Range("I37").Value = "4#5#6#7#8#1#14#15#16#17#18#21#22#23#24#25#26#27#28#29#30#31#32#33#34"
$ws.Range("J37").Value = '"+#+#+#+#-#-#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#-'
$ws.Range("K37").Value = 6

# --- Row 41 (LABA USAHA SETELAH PAJAK): new formula/operator, Report Type bumped ---
$ws.Range("I41").Value = "4#5#6#7#8#1#14#15#16#17#18#21#22#23#24#25#26#27#28#29#30#31#32#33#34#38"
$ws.Range("J41").Value = '"+#+#+#+#-#-#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#+#-#-'
$ws.Range("K41").Value = 6

# --- Row heights: bump data rows to 16.5pt (rows 1, 3, 8, 12 and 20 stay default) ---
$ws.Range("2:2").RowHeight = 16.5
$ws.Range("4:7").RowHeight = 16.5
$ws.Range("9:11").RowHeight = 16.5
$ws.Range("13:19").RowHeight = 16.5
$ws.Range("21:41").RowHeight = 16.5

# --- Widen the "Report Operator" column (J) so the longer strings are visible ---
$ws.Range("J1").ColumnWidth = 45.8

# --- Move the active selection (previously scrolled to G29 / K35 selected) ---
$ws.Range("I12").Select()
